$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 59, shifting rows 59-61 down to 60-62.
$ws.Rows.Item(59).Insert()

# The SUM formula's range did not auto-extend to include the newly
# inserted blank row, so widen it explicitly to match row 59 (new).
$ws.Range("F60").Formula = "=SUM(F2:F59)"

# Select A61 as the active cell (matches the post-edit selection in the file)
$ws.Range("A61").Select()
